{"js": "const pairs = [\n  [\"909\u00d74=\", \"544\u00d78=\"],\n  [\"626\u00d75=\", \"863\u00d76=\"],\n  [\"690\u00d72=\", \"438\u00d75=\"],\n  [\"435\u00d72=\", \"768\u00d76=\"],\n  [\"392\u00d77=\", \"330\u00d78=\"],\n  [\"482\u00d72=\", \"512\u00d72=\"],\n  [\"825\u00d77=\", \"754\u00d79=\"],\n  [\"151\u00d73=\", \"870\u00d76=\"],\n  [\"832\u00d76=\", \"205\u00d74=\"],\n  [\"563\u00d76=\", \"281\u00d73=\"],\n  [\"726\u00d74=\", \"825\u00d78=\"],\n  [\"470\u00d77=\", \"224\u00d72=\"],\n  [\"674\u00d73=\", \"595\u00d75=\"],\n  [\"441\u00d72=\", \"302\u00d76=\"],\n  [\"341\u00d79=\", \"982\u00d76=\"],\n  [\"931\u00d79=\", \"332\u00d79=\"],\n  [\"358\u00d78=\", \"257\u00d75=\"],\n  [\"797\u00d79=\", \"352\u00d77=\"],\n  [\"743\u00d79=\", \"514\u00d73=\"],\n  [\"217\u00d76=\", \"878\u00d72=\"],\n  [\"427\u00d75=\", \"905\u00d76=\"],\n  [\"420\u00d77=\", \"591\u00d79=\"],\n  [\"247\u00d78=\", \"932\u00d75=\"],\n  [\"269\u00d73=\", \"285\u00d77=\"],\n  [\"810\u00d75=\", \"340\u00d74=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"909\u00d74=\", \"544\u00d78=\"),\n    @(\"626\u00d75=\", \"863\u00d76=\"),\n    @(\"690\u00d72=\", \"438\u00d75=\"),\n    @(\"435\u00d72=\", \"768\u00d76=\"),\n    @(\"392\u00d77=\", \"330\u00d78=\"),\n    @(\"482\u00d72=\", \"512\u00d72=\"),\n    @(\"825\u00d77=\", \"754\u00d79=\"),\n    @(\"151\u00d73=\", \"870\u00d76=\"),\n    @(\"832\u00d76=\", \"205\u00d74=\"),\n    @(\"563\u00d76=\", \"281\u00d73=\"),\n    @(\"726\u00d74=\", \"825\u00d78=\"),\n    @(\"470\u00d77=\", \"224\u00d72=\"),\n    @(\"674\u00d73=\", \"595\u00d75=\"),\n    @(\"441\u00d72=\", \"302\u00d76=\"),\n    @(\"341\u00d79=\", \"982\u00d76=\"),\n    @(\"931\u00d79=\", \"332\u00d79=\"),\n    @(\"358\u00d78=\", \"257\u00d75=\"),\n    @(\"797\u00d79=\", \"352\u00d77=\"),\n    @(\"743\u00d79=\", \"514\u00d73=\"),\n    @(\"217\u00d76=\", \"878\u00d72=\"),\n    @(\"427\u00d75=\", \"905\u00d76=\"),\n    @(\"420\u00d77=\", \"591\u00d79=\"),\n    @(\"247\u00d78=\", \"932\u00d75=\"),\n    @(\"269\u00d73=\", \"285\u00d77=\"),\n    @(\"810\u00d75=\", \"340\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $old,\n        $false,\n        $true,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $new,\n        2\n    ) | Out-Null\n}\n"}
